$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actions")

# Copy the formatting of the row above (row 4) into the new row 5 cells for
# the columns whose style changes, then fill in the new scenario data.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E4").Copy()
$ws.Range("E5").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A5").Value = "sc4"
$ws.Range("B5").Value = "navigateToYMORegistrationPage"
$ws.Range("C5").Value = "yes"
$ws.Range("D5").Value = "fail"
$ws.Range("E5").Value = "YMO"

# New, longer scenario name needs the column to be widened to fit it
# (best-fit to roughly 28.7 characters wide).
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(2).ColumnWidth = 27.8

# The validation list should only apply to column C now (not C:D).
$ws.Range("C2:D31").Validation.Delete()
$ws.Range("C2:C31").Validation.Add(3, 1, 1, '"yes,no"')

# Make "Actions" the active sheet/tab, with D9 selected.
$ws.Activate()
$ws.Range("D9").Select()
